$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 9.458811508128221
$ws.Range("R2").Value = 85.129303573154
$ws.Range("S2").Value = 0.08142722252275085
$ws.Range("T2").Value = 0.08142722252275086
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("S3").Value = 0.07895519991170555
$ws.Range("T3").Value = 0.07895519991170555
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 9.339736627803333
$ws.Range("R4").Value = 84.05762965023
$ws.Range("S4").Value = 0.08040215327713235
$ws.Range("T4").Value = 0.08040215327713235
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 3.726291358421223
$ws.Range("R5").Value = 33.536622225791
$ws.Range("S5").Value = 0.03207819030604742
$ws.Range("T5").Value = 0.03207819030604742
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 25.20627543714178
$ws.Range("R6").Value = 226.856478934276
$ws.Range("S6").Value = 0.2169910032805011
$ws.Range("T6").Value = 0.2169910032805011
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.2104034438638339
$ws.Range("T7").Value = 0.2104034438638339
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 24.88895922584667
$ws.Range("R8").Value = 224.00063303262
$ws.Range("S8").Value = 0.2142593516639106
$ws.Range("T8").Value = 0.2142593516639106
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 9.929992394783779
$ws.Range("R9").Value = 89.36993155305402
$ws.Range("S9").Value = 0.08548343517411824
$ws.Range("T9").Value = 0.08548343517411824
